# "Atualização de horas trabalhadas no desenvolvimento"
# Update the hours-logged-per-day table on the "Dados" sheet: move some of
# the logged hours for each activity to later days in the sprint.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados")

# Row 2 - "Cadastrar Usuário" (Horas = 5): log 2h on day 18 (col U) and
# 3h on day 38 (col AO) instead of being unlogged.
$ws.Cells.Item(2, 21).Value = 2    # U2
$ws.Cells.Item(2, 41).Value = 3    # AO2

# Row 3 - "Realizar Login" (Horas = 8): move the 1h previously on day 1
# (col D) to 4h on day 19 (col V) and 4h on day 39 (col AP).
$ws.Cells.Item(3, 4).Value = $null # D3 (clear)
$ws.Cells.Item(3, 22).Value = 4    # V3
$ws.Cells.Item(3, 42).Value = 4    # AP3

# Row 4 - "Candidatar a Monitor" (Horas = 10): move the 1h previously on
# day 9 (col L) to 5h on day 20 (col W), and bump day 40 (col AQ) from
# 2h to 5h.
$ws.Cells.Item(4, 12).Value = $null # L4 (clear)
$ws.Cells.Item(4, 23).Value = 5     # W4
$ws.Cells.Item(4, 43).Value = 5     # AQ4 (was 2)

# Refresh the sheet view on "Dados": zoom to 85% and move the selection.
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("V21:W21").Select()

# Restore "Grafico" as the active/selected sheet tab.
$wsGrafico = $wb.Worksheets.Item("Grafico")
$wsGrafico.Activate()
